# Applies the edits described by the commit:
#  - products: update image column (D) to full file paths (except "controller" row)
#              and widen column D
#  - orders_details: drop the leftover empty "image" column (D)
#  - DDict_attributes: fix header D1 "description" -> "caption"
#  - add a new "extra_test" worksheet at the end of the workbook
#  - restore selections / the originally active sheet ("products")

$wb = $excel.ActiveWorkbook

# --- products: rewrite the image paths (full Windows paths from the data folder) ---
$products = $wb.Worksheets.Item("products")
$imgBase = "C:\Users\SSAran\Desktop\Files\duty\Uniss\EnhanceFAIRness_spreadsheet_to_db\data\images\"
$products.Range("D3").Value = $imgBase + "nail file.jpg"
$products.Range("D4").Value = $imgBase + "backpack.jpg"
$products.Range("D5").Value = $imgBase + "pen.jpg"
$products.Range("D6").Value = $imgBase + "phone.jpg"
$products.Range("D7").Value = $imgBase + "sunglasses.jpg"
$products.Columns.Item(4).ColumnWidth = 114.877604166667

# --- orders_details: the "image" column (D) was a stray leftover, remove it ---
$ordersDetails = $wb.Worksheets.Item("orders_details")
$ordersDetails.Columns.Item(4).Delete() | Out-Null

# --- DDict_attributes: fix header caption ---
$ddictAttr = $wb.Worksheets.Item("DDict_attributes")
$ddictAttr.Range("D1").Value = "caption"

# --- add the new "extra_test" worksheet at the end ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$extraTest = $wb.Worksheets.Add($null, $lastSheet)
$extraTest.Name = "extra_test"
$extraTest.Range("B1").Value = "column_1"
$extraTest.Range("C1").Value = "column_2"
$extraTest.Range("B2").Value = "a table that should not be keeped"
$extraTest.Range("C2").Value = 42
$extraTest.Columns.Item(2).ColumnWidth = 31.736979166667
$extraTest.Range("F4").Select() | Out-Null

# --- restore selections on the touched sheets ---
$ordersDetails.Range("E9").Select() | Out-Null
$ddictAttr.Range("D1").Select() | Out-Null

$products.Activate() | Out-Null
$products.Range("D14").Select() | Out-Null

Write-Output "done"
